# Apply updated coin price/volume/name/link data per the Tue Feb 14 06:44:51 UTC 2023 symbol-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.27"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-6.90%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'40.40"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-0.96%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.030"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-2.42%"
$ws.Range("E4").ClearFormats()
$ws.Range("E5").Value = "'-3.57%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.529"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-9.00%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.9298"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.07%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'2.384"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-1.66%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1173"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-2.13%"
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'-3.83%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.04331"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'4.58%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.08683"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-4.08%"
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'0.08%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001279"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.03942"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-2.14%"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005921"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'1.54%"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.335"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.14%"
$ws.Range("E17").ClearFormats()
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.280"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-1.09%"
$ws.Range("E18").ClearFormats()
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3288"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-2.06%"
$ws.Range("E19").ClearFormats()
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.975"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'4.74%"
$ws.Range("E20").ClearFormats()
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1400"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'4.33%"
$ws.Range("E21").ClearFormats()
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2742"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-3.36%"
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'-1.53%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.003783"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-4.81%"
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'-1.57%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.0003725"
$ws.Range("D26").ClearFormats()
$ws.Range("D38").Value = "'0.02282"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'-5.37%"
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.05041"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").Value = "'0.005856"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'77.34%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007677"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.90%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'-1.07%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.007356"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-3.32%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.008258"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-3.79%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.2919"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-13.81%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006269"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-4.90%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'0.02%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.03188"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'-88.13%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'0.02%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E50").ClearFormats()
